$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "杯中的风味世界"; New = "柴茶：杯子里的味道世界" },
    @{ Old = "健康与愉悦的完美融合"; New = "柴茶：健康与快乐的完美融合" },
    @{ Old = "不仅仅是茶，更是一种生活方式"; New = "柴茶：不仅仅是茶，一种生活方式" },
    @{ Old = "四季皆宜的饮品"; New = "柴茶：所有季节和原因的饮料" },
    @{ Old = "感官的极致享受"; New = "柴茶：你的感官的终极放纵" },
    @{ Old = "远离日常生活的甜蜜之选"; New = "柴茶：从日常的甜蜜逃跑" },
    @{ Old = "分享温暖，分享爱"; New = "柴茶：分享温暖，分享爱" }
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
